$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.144224824685715
$ws.Range("A3").Value = 0.4179271137421
$ws.Range("A4").Value = 0.46298807126105
$ws.Range("A5").Value = 0.0572041752120773

$ws.Range("A9").Value = 0.0500923428320895
$ws.Range("A10").Value = 0.253051866997579
$ws.Range("A11").Value = 0.688350595617186
$ws.Range("A12").Value = 0.0385357915393103

$ws.Range("A13").Value = 0.401534343478692
$ws.Range("A14").Value = 0.0289684862356343
$ws.Range("A15").Value = 0.201023573873398
$ws.Range("A16").Value = 0.338706610125542
$ws.Range("A17").Value = 0.197418414026358

$ws.Range("A18").Value = 0.0365590276042801
$ws.Range("A19").Value = 0.0352524404264423
$ws.Range("A20").Value = 0.620173540603422
$ws.Range("A21").Value = 0.345783685446419
$ws.Range("A22").Value = 0.0409233601077701

$ws.Range("A23").Value = 0.686279141430231
$ws.Range("A24").Value = 0.27894878156222
$ws.Range("A25").Value = 0.0694555525133838
$ws.Range("A26").Value = -0.045332574500825

$ws.Range("A36").Value = 0.651074238362871
$ws.Range("A37").Value = 0.286518821180841
$ws.Range("A38").Value = 0.0859439287188909
$ws.Range("A39").Value = 0.0484593984171848
$ws.Range("A40").Value = 0.0348745835085473

$ws.Range("A41").Value = 0.537210491859705
$ws.Range("A42").Value = 0.267150565667853
$ws.Range("A43").Value = 0.0937444877361151
$ws.Range("A44").Value = 0.206599104763492
$ws.Range("A45").Value = 0.0493247308602022

$ws.Range("A46").Value = 0.184421616982797
$ws.Range("A47").Value = 0.0916295543816123
$ws.Range("A48").Value = 0.11790509323629
$ws.Range("A49").Value = 0.225317050602373
$ws.Range("A50").Value = 0.527175043370673

$ws.Range("A51").Value = 0.0504090655570513
$ws.Range("A52").Value = 0.559218958872528
$ws.Range("A53").Value = 0.263543404292093
$ws.Range("A54").Value = 0.10162766835581
$ws.Range("A55").Value = 0.0511846159308557
$ws.Range("A56").Value = 0.133184394120326
